$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# Text-like numeric values (e.g. "244.07") must remain stored as text, matching
# the source data which uses inline/shared strings rather than numeric cells.
# We briefly force a Text number format while assigning the value so Excel does not
# auto-convert the string to a floating point number, then clear the format again
# so the cell keeps its original (default/general) style with no explicit style index.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "36.415.64"
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -2.65%  "
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.981.71"
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -3.42%  "
$cell.ClearFormats()

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "244.07"
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -3.26%  "
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -3.69%  "
$cell.ClearFormats()

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "59.08"
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -9.66%  "
$cell.ClearFormats()

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.00%  "
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -2.17%  "
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "57.49"
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -3.18%  "
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0817"
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +6.51%  "
$cell.ClearFormats()

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -0.97%  "
$cell.ClearFormats()

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "23.81"
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +5.35%  "
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.862"
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -6.31%  "
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "13.98"
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -5.72%  "
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.272.49"
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -3.56%  "
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.46"
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -2.16%  "
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.979.73"
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -3.91%  "
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "36.300.78"
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -2.77%  "
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "69.89"
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -4.92%  "
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -1.22%  "
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.31"
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -3.08%  "
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "234.46"
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.23%  "
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -2.11%  "
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.30"
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -4.11%  "
$cell.ClearFormats()

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.04"
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.16%  "
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "162.16"
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +0.24%  "
$cell.ClearFormats()

$cell = $ws.Range("B29")
$cell.NumberFormat = "@"
$cell.Value = "Kaspa"
$cell.ClearFormats()

$cell = $ws.Range("C29")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell.ClearFormats()

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.133"
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +11.15%  "
$cell.ClearFormats()

$cell = $ws.Range("B30")
$cell.NumberFormat = "@"
$cell.Value = "EthereumClassic"
$cell.ClearFormats()

$cell = $ws.Range("C30")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell.ClearFormats()

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "19.82"
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.76%  "
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.91"
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -6.73%  "
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0628"
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -5.34%  "
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.21"
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +2.92%  "
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -7.45%  "
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.79"
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -2.45%  "
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.03"
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -7.65%  "
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -4.05%  "
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -2.30%  "
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -4.89%  "
$cell.ClearFormats()

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "92.51"
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -4.20%  "
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "16.19"
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -6.65%  "
$cell.ClearFormats()

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.371.08"
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -3.43%  "
$cell.ClearFormats()

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.49"
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -5.57%  "
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.84"
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -3.56%  "
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -2.75%  "
$cell.ClearFormats()
